# Weekly fruit/vegetable price update:
# Insert a new record (row 86) for "Poroto verde" at Feria Lagunitas de
# Puerto Montt, pushing the existing rows 86-106 down to 87-107.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(86).Insert()

$ws.Cells.Item(86, 1).Value = 4
$ws.Cells.Item(86, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(86, 3).Value = "Los Lagos"
$ws.Cells.Item(86, 4).Value = 44855
$ws.Cells.Item(86, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(86, 5).Value = 10
$ws.Cells.Item(86, 6).Value = 100112031
$ws.Cells.Item(86, 7).Value = "Poroto verde"
$ws.Cells.Item(86, 8).Value = "Magnum"
$ws.Cells.Item(86, 9).Value = "Primera"
$ws.Cells.Item(86, 10).Value = 35
$ws.Cells.Item(86, 11).Value = 35000
$ws.Cells.Item(86, 12).Value = 35000
$ws.Cells.Item(86, 13).Value = 35000
$ws.Cells.Item(86, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(86, 15).Value = "Perú"
$ws.Cells.Item(86, 16).Value = 1400
$ws.Cells.Item(86, 17).Value = 25
$ws.Cells.Item(86, 18).Value = "Hortaliza"
